$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet/tab to reflect the new "through" date
$ws.Name = "Through 2022-06-25"

# Update the row label for June to reflect the new "through" date
$ws.Range("A7").Value = "June (through 06-25)"

# Update the June row (row 7) values that changed
$ws.Range("D7").Value = 60
$ws.Range("E7").Value = 47
$ws.Range("G7").Value = 93
$ws.Range("H7").Value = 101
$ws.Range("I7").Value = 118

# Update the Total row (row 8) values that changed
$ws.Range("D8").Value = 376
$ws.Range("E8").Value = 342
$ws.Range("G8").Value = 451
$ws.Range("H8").Value = 732
$ws.Range("I8").Value = 781
